$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update CMND number, turn "Ngay sinh" into a real date value
$ws.Range("A2").Value = 987654321
$ws.Range("C2").Value = 36454
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Row 3: update CMND number, replace the placeholder name, turn "Ngay sinh" into a real date value
$ws.Range("A3").Value = 987654322
$ws.Range("B3").Value = "Trương Văn Khôi"
$ws.Range("C3").Value = 36475

# Copy C2's formatting onto C3 so both date cells share one cell style
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Widen column C to fit the new date content
$ws.Columns.Item(3).ColumnWidth = 11.4

# Move the active selection to C4, matching the post-edit cursor position
$ws.Range("C4").Select()
